$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before column B ("Seed" -> C, "Independent layer 0" -> D, "Incremental layer 0" -> E)
$ws.Columns("B:B").Insert()

# 2. New header "batch size" in B2
$ws.Range("B2").Value = "batch size"

# 3. Fill the new "batch size" column for the existing rows (3-6) with 128
$ws.Range("B3").Value = 128
$ws.Range("B4").Value = 128
$ws.Range("B5").Value = 128
$ws.Range("B6").Value = 128

# 4. New row 7: another CodeGPTPy / batch 1024 / seed 0 entry
$ws.Range("A7").Value = "CodeGPTPy"
$ws.Range("B7").Value = 1024
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0.685
$ws.Range("E7").Value = 0.70499999999999996

# Copy the row5 formatting (bordered + black font) onto row 7's A:C cells,
# matching the "CodeGPTPy"/seed formatting used for the other 1024-batch row.
$ws.Range("A5").Copy()
$ws.Range("A7").PasteSpecial(-4122)

# 5. New "Diff." column (F) = Independent - Incremental accuracy
$ws.Range("F2").Value = "Diff."
$ws.Range("F3").Formula = "=D3-E3"
$ws.Range("F4:F7").Formula = "=D4-E4"

# 6. Number format (0.000) for the new Diff. column
$ws.Range("F3:F7").NumberFormat = "0.000"

# 7. Give the Diff. header + the new row's batch/seed cells a left/right-only
#    thin border, and the new row's accuracy cells the same border plus the
#    0.000 number format, by copying from a fully-bordered cell and then
#    stripping the top/bottom edges off the copy.
$ws.Range("A2").Copy()
$ws.Range("F2").PasteSpecial(-4122)
$ws.Range("F2").Borders.Item(8).LineStyle = -4142
$ws.Range("F2").Borders.Item(9).LineStyle = -4142
$ws.Range("F2").Value = "Diff."

$ws.Range("B7:C7").Value2 = $ws.Range("B7:C7").Value2
$ws.Range("A2").Copy()
$ws.Range("B7:C7").PasteSpecial(-4122)
$ws.Range("B7:C7").Borders.Item(8).LineStyle = -4142
$ws.Range("B7:C7").Borders.Item(9).LineStyle = -4142
$ws.Range("B7").Value = 1024
$ws.Range("C7").Value = 0

$ws.Range("D3").Copy()
$ws.Range("D7:E7").PasteSpecial(-4122)
$ws.Range("D7:E7").Borders.Item(8).LineStyle = -4142
$ws.Range("D7:E7").Borders.Item(9).LineStyle = -4142
$ws.Range("D7").Value = 0.685
$ws.Range("E7").Value = 0.70499999999999996

$ws.Range("F3:F7").Borders.Item(7).LineStyle = -4142
$ws.Range("F3:F7").Borders.Item(10).LineStyle = -4142
$ws.Range("F3:F7").Borders.Item(8).LineStyle = -4142
$ws.Range("F3:F7").Borders.Item(9).LineStyle = -4142

$wb.Save()
